# Add the "2016_monthly" worksheet (new report tab) with headers + 12 months
# of resident / non-resident termination counts, matching the commit:
# "add monthly data to reports ... add monthly visualizations and
#  datatable to dashboard ..."

$wb = $excel.ActiveWorkbook

# Remember the sheet that currently holds the tab selection so we can move
# its "last selected cell" off the header block once the new tab becomes
# active (mirrors what Excel itself records when you switch tabs).
$prevActive = $excel.ActiveSheet

# Insert the new sheet as the LAST tab (after "2016_county") rather than
# Excel's default "before the active sheet" placement.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2016_monthly"

# Header row
$ws.Range("A1").Value = "Month"
$ws.Range("B1").Value = "Resident Terminations`r(n = 6,767)"
$ws.Range("C1").Value = "Non-Resident Terminations`r(n = 510)"

# Monthly data: Month, Resident Terminations, Non-Resident Terminations
$monthlyData = @(
    @("January",   599, 34),
    @("February",  607, 49),
    @("March",     687, 54),
    @("April",     661, 35),
    @("May",       566, 55),
    @("June",      552, 61),
    @("July",      549, 52),
    @("August",    530, 28),
    @("September", 530, 49),
    @("October",   476, 33),
    @("November",  502, 18),
    @("December",  508, 42)
)

$r = 2
foreach ($row in $monthlyData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Make the new report tab the active one (Excel moves tabSelected / the
# workbook's activeTab pointer to whichever sheet is active on save, and
# the previously-active sheet keeps whatever cell was last selected there).
# NB: selecting a range on another sheet implicitly activates it, so set
# the outgoing sheet's remembered selection BEFORE activating the new tab.
if ($prevActive -ne $null -and $prevActive.Name -ne $ws.Name) {
    $prevActive.Range("H19").Select()
}

$ws.Activate()
$ws.Range("E5").Select()

Write-Output "Added sheet '$($ws.Name)'; workbook now has $($wb.Worksheets.Count) sheets."
